# 389-RBI-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-EarlyRePayment-Newcreateloan.xlsx
# "Loan RBI, Variable Instalments"
#
# On the "Repayment Schedule" sheet a new (blank) instalment-related column is
# inserted before column N ("Late"), pushing "Late" to column O and
# "Outstanding" to column Q (leaving a blank column N and a blank column P).
# The "Repayment Schedule" sheet also becomes the active/selected sheet
# (previously "NewLoanInput" was selected), with L17 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Make "Repayment Schedule" the active sheet (this also clears tabSelected on
# whichever sheet previously had it, i.e. "NewLoanInput").
$ws.Activate() | Out-Null

# Insert a new blank column before column N (14) - shifts "Late"/"Outstanding"
# (and the data beneath them) one column to the right.
$ws.Columns.Item(14).Insert() | Out-Null

# Update the selection on the sheet to match the new active cell.
$ws.Range("L17").Select() | Out-Null
